$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: PT Mandiri Utama Finance - update lease dates, actual income, and payment scheme
$ws.Range("B15").Value = 45991
$ws.Range("C15").Value = 47817
$ws.Range("G15").Value = 280000000
$ws.Range("H15").Value = "Full Lease Upfront"
$ws.Range("I15").Value = ""
